# MRD-470 Response to probation
#
# The "Provide details of how the offender has responded to supervision to
# date?" question is a legacy FORMTEXT field whose (empty) result is stored
# as 5 runs of an en-space character (U+2002). Replace that result with the
# {{response_to_probation}} merge placeholder.

$d = $word.ActiveDocument

# Locate the question heading so we can scope the subsequent replacement to
# just this field's result (the same "empty result" run pattern recurs many
# times across the document for other questions).
$heading = $d.Content
$found = $heading.Find.Execute("Provide details of how the offender has responded to supervision to date?")
if (-not $found) {
    throw "Could not find the 'responded to supervision to date' question heading"
}

# Move to the field immediately following the heading paragraph and grab a
# generous window that comfortably contains the field's begin/separate/end
# codes and its (normally blank) result runs.
$heading.Collapse(0)
[void]$heading.MoveEnd(1, 40)

$enSpace = [char]0x2002
$emptyResult = "$enSpace$enSpace$enSpace$enSpace$enSpace"

$replaced = $heading.Find.Execute($emptyResult, $false, $false, $false, $false, $false, $true, 1, $false, "{{response_to_probation}}", 2)
if (-not $replaced) {
    throw "Could not find the response-to-probation field result to replace"
}
